# "add user list to project"
#
# 1) The cfop list for Romit Roy Choudhury (PI hours!G3) had its two
#    entries reordered.
# 2) The "project hours" sheet gains a new "users" column (E) that lists
#    the user(s) associated with each project row.

$wb = $excel.ActiveWorkbook

# --- 1. Fix the reordered cfop list on the "PI hours" sheet -----------------
$piSheet = $wb.Worksheets.Item("PI hours")
$piSheet.Range("G3").Value = "['cfop_RRC', 'cfop_CHOUDHURY']"

# --- 2. Add the new "users" column to the "project hours" sheet -------------
$projSheet = $wb.Worksheets.Item("project hours")

# Give the new header cell the same look (bold / border / centered) as the
# other header cells in row 1 by copying the format from D1 first.
$projSheet.Range("D1").Copy()
$projSheet.Range("E1").PasteSpecial(-4122)
$projSheet.Range("E1").Value = "users"

$projSheet.Range("E2").Value = "['Arun Lakshmanan', 'Mitchell Jones']"
$projSheet.Range("E3").Value = "['Ashutosh Dhekne']"
$projSheet.Range("E4").Value = "['Chawla, Karan Gansham']"
$projSheet.Range("E5").Value = "['Won Dong Shin']"
$projSheet.Range("E6").Value = "['Aygen Berk Cagilci']"
